# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -5
$ws.Range("F15").Value = 3
